$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D5").Value = "NCAP_ILED"
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = "Trans - Update"
